$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.094.69'
$ws.Range('E2').Value = '  +0.04%  '
$ws.Range('D3').Value = '1.831.31'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.9986'
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '242.77'
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.6283'
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.9995'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('E8').Value = '  -1.04%  '
$ws.Range('E9').Value = '  -0.27%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '23.24'
$ws.Range('E10').Value = '  +2.58%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07679'
$ws.Range('E11').Value = '  -0.85%  '
$ws.Range('D12').Value = '1.831.82'
$ws.Range('E12').Value = '  -0.28%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.025'
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.6694'
$ws.Range('E14').Value = '  +0.47%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '82.96'
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.000009390'
$ws.Range('E16').Value = '  -6.82%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.994'
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('D18').Value = '29.097.53'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').Value = '2.076.59'
$ws.Range('E19').Value = '  -0.79%  '
$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '223.21'
$ws.Range('E21').Value = '  -2.00%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.0000'
$ws.Range('E22').Value = '  -0.25%  '
$ws.Range('E23').Value = '  -1.17%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.9995'
$ws.Range('E24').Value = '  -0.24%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '159.88'
$ws.Range('E25').Value = '  +0.28%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.1399'
$ws.Range('E26').Value = '  +1.03%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.495'
$ws.Range('E27').Value = '  -0.37%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '17.88'
$ws.Range('E28').Value = '  -0.43%  '
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05753'
$ws.Range('E30').Value = '  +9.29%  '
$ws.Range('E31').Value = '  +1.17%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.117'
$ws.Range('E32').Value = '  +2.24%  '
$ws.Range('E33').Value = '  +0.59%  '
$ws.Range('B34').Value = 'LidoDAOToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.837'
$ws.Range('E34').Value = '  -0.41%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.7415'
$ws.Range('E35').Value = '  +0.76%  '
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.664'
$ws.Range('E37').Value = '  -1.26%  '
$ws.Range('D38').Value = '1.229.21'
$ws.Range('E38').Value = '  -1.21%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.763'
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.01779'
$ws.Range('E40').Value = '  -0.58%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.501'
$ws.Range('E41').Value = '  +2.03%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.8916'
$ws.Range('E42').Value = '  -0.75%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.9989'
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '101.87'
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '65.79'
$ws.Range('E46').Value = '  +2.00%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.00000000124'
$ws.Range('E47').Value = '  +0.51%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.5088'
$ws.Range('E48').Value = '  -0.58%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.07583'
$ws.Range('E49').Value = '  +13.71%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.4067'
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '9.016'
$ws.Range('E51').Value = '  +0.95%  '
